# Adding new link to scraping
# Appends 6 new rows (234-239) of COVID data below the existing data
# (dates descending continue from 02.04.2020 on the old last row down to
# 03.01.2021 on the new last row) and updates the sheet dimension.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data to append: row number, date (dd.mm.yyyy text), then 6 numeric values
$newRows = @(
    @(234, "11.01.2021", 2875, 11529, 144, 1963, 768, 0),
    @(235, "10.01.2021", 2869, 11486, 141, 1852, 876, 0),
    @(236, "09.01.2021", 2854, 11408, 140, 1750, 964, 0),
    @(237, "08.01.2021", 2834, 11345, 139, 1648, 1047, 0),
    @(238, "05.01.2021", 2780, 11177, 136, 1494, 1150, 0),
    @(239, "03.01.2021", 2737, 11026, 134, 1475, 1128, 0)
)

foreach ($row in $newRows) {
    $r = $row[0]

    # Column A holds the date as literal text (e.g. "11.01.2021"), not a
    # real date value, matching the rest of the column. Forcing the cell
    # to Text format before assignment keeps Excel from auto-converting
    # the string into a date serial number; ClearFormats afterwards drops
    # the now-unneeded explicit format so the cell stays unstyled, just
    # like its neighbours.
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $row[1]
    $dateCell.ClearFormats()

    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}
